{"js": "// Replace the 25 two-digit-by-two-digit multiplication prompts in the\n// table with their new operands. Each \"before\" string below is unique\n// within the document, so a matchCase, non-wildcard search locates\n// exactly one run per pair; insertText(..., \"Replace\") swaps only the\n// text (the run's formatting, e.g. TimeNewRoman/sz 30, is preserved).\nconst replacements = [\n  [\"23\u00d718=\", \"70\u00d762=\"],\n  [\"38\u00d728=\", \"52\u00d751=\"],\n  [\"56\u00d795=\", \"37\u00d787=\"],\n  [\"64\u00d744=\", \"31\u00d786=\"],\n  [\"96\u00d797=\", \"36\u00d723=\"],\n  [\"27\u00d751=\", \"49\u00d752=\"],\n  [\"66\u00d711=\", \"84\u00d788=\"],\n  [\"94\u00d767=\", \"17\u00d714=\"],\n  [\"82\u00d720=\", \"98\u00d762=\"],\n  [\"88\u00d794=\", \"13\u00d723=\"],\n  [\"57\u00d727=\", \"42\u00d712=\"],\n  [\"63\u00d720=\", \"56\u00d717=\"],\n  [\"27\u00d790=\", \"19\u00d732=\"],\n  [\"65\u00d738=\", \"76\u00d736=\"],\n  [\"66\u00d756=\", \"38\u00d739=\"],\n  [\"86\u00d725=\", \"38\u00d749=\"],\n  [\"55\u00d736=\", \"80\u00d723=\"],\n  [\"92\u00d731=\", \"82\u00d763=\"],\n  [\"28\u00d716=\", \"31\u00d711=\"],\n  [\"52\u00d724=\", \"91\u00d734=\"],\n  [\"76\u00d784=\", \"76\u00d743=\"],\n  [\"47\u00d786=\", \"70\u00d721=\"],\n  [\"85\u00d755=\", \"60\u00d718=\"],\n  [\"41\u00d714=\", \"11\u00d769=\"],\n  [\"34\u00d765=\", \"51\u00d726=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit-by-two-digit multiplication prompts in the\n# table with their new operands. Each \"before\" string is unique within\n# the document, so Find/Replace (MatchCase, no wildcards) on the whole\n# document body touches exactly one run per pair and only swaps the\n# text -- the run's formatting (TimeNewRoman / sz 30) is left alone.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"23\u00d718=\", \"70\u00d762=\"),\n  @(\"38\u00d728=\", \"52\u00d751=\"),\n  @(\"56\u00d795=\", \"37\u00d787=\"),\n  @(\"64\u00d744=\", \"31\u00d786=\"),\n  @(\"96\u00d797=\", \"36\u00d723=\"),\n  @(\"27\u00d751=\", \"49\u00d752=\"),\n  @(\"66\u00d711=\", \"84\u00d788=\"),\n  @(\"94\u00d767=\", \"17\u00d714=\"),\n  @(\"82\u00d720=\", \"98\u00d762=\"),\n  @(\"88\u00d794=\", \"13\u00d723=\"),\n  @(\"57\u00d727=\", \"42\u00d712=\"),\n  @(\"63\u00d720=\", \"56\u00d717=\"),\n  @(\"27\u00d790=\", \"19\u00d732=\"),\n  @(\"65\u00d738=\", \"76\u00d736=\"),\n  @(\"66\u00d756=\", \"38\u00d739=\"),\n  @(\"86\u00d725=\", \"38\u00d749=\"),\n  @(\"55\u00d736=\", \"80\u00d723=\"),\n  @(\"92\u00d731=\", \"82\u00d763=\"),\n  @(\"28\u00d716=\", \"31\u00d711=\"),\n  @(\"52\u00d724=\", \"91\u00d734=\"),\n  @(\"76\u00d784=\", \"76\u00d743=\"),\n  @(\"47\u00d786=\", \"70\u00d721=\"),\n  @(\"85\u00d755=\", \"60\u00d718=\"),\n  @(\"41\u00d714=\", \"11\u00d769=\"),\n  @(\"34\u00d765=\", \"51\u00d726=\")\n)\n\n$wdReplaceAll = 2\n$wdFindWrap = 1\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindWrap, $false, $null, $wdReplaceAll) | Out-Null\n}\n"}
